# edit.ps1 - reproduces the author's edit:
#   1) Re-style the three data tables (slides 14, 15, 16) from the
#      locally-defined "Table_0" style to the built-in table style
#      {C1B19D57-1F49-4794-9282-DCBEDADDFAA7}.
#   2) Switch the presentation's applied Design theme's colour scheme
#      from the "Integral"/"Red Violet" palette to the stock
#      "Office Theme"/"Office" palette (Design gallery: Office Theme).

$p = $ppt.ActivePresentation

# --- 1) Table style swap -------------------------------------------------
$newStyleId = "{C1B19D57-1F49-4794-9282-DCBEDADDFAA7}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Theme colour scheme swap ----------------------------------------
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$slide1 = $p.Slides.Item(1)
$colors = $slide1.ColorScheme

$colors.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$colors.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$colors.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$colors.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$colors.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$colors.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$colors.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$colors.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$colors.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hyperlink
$colors.Colors(12).RGB = RGB 0x95 0x4F 0x72   # followed hyperlink
